$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1561.5
$ws.Range("J17").Value = 1561.5
$ws.Range("L17").Value = 4684.5
$ws.Range("N17").Value = -5020.5
$ws.Range("H112").Value = 2123.6785
$ws.Range("J112").Value = 2221.6538
$ws.Range("L112").Value = 6664.9614
$ws.Range("N112").Value = -8880.9614
$ws.Range("H129").Value = 1095.1875
$ws.Range("J129").Value = 1203.3846
$ws.Range("L129").Value = 3610.1538
$ws.Range("N129").Value = -13610.1538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 613.4888999999999
$ws.Range("J2").Value = 1519.8182
$ws.Range("L2").Value = 1519.8182
$ws.Range("N2").Value = -1745.8182
$ws.Range("H32").Value = 19341.441
$ws.Range("I32").Value = 18611.477
$ws.Range("K32").Value = 18611.477
$ws.Range("M32").Value = -18324.477
$ws.Range("H110").Value = 1395.25
$ws.Range("I110").Value = 1355.6923
$ws.Range("J110").Value = 1566.6666
$ws.Range("K110").Value = 1355.6923
$ws.Range("L110").Value = 1566.6666
$ws.Range("M110").Value = 689.3077000000001
$ws.Range("N110").Value = -5656.6666
$ws.Range("H116").Value = 613.4888999999999
$ws.Range("J116").Value = 1519.8182
$ws.Range("L116").Value = 1519.8182
$ws.Range("N116").Value = -6107.8182
$ws.Range("H122").Value = 1781.5834
$ws.Range("I122").Value = 1860.3158
$ws.Range("K122").Value = 5580.9474
$ws.Range("M122").Value = -3130.9474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 613.4888999999999
$ws.Range("J3").Value = 1519.8182
$ws.Range("L3").Value = 1519.8182
$ws.Range("N3").Value = -1747.8182
$ws.Range("H94").Value = 1296.826
$ws.Range("I94").Value = 1053.9412
$ws.Range("J94").Value = 1985
$ws.Range("K94").Value = 1053.9412
$ws.Range("L94").Value = 1985
$ws.Range("M94").Value = -602.9412
$ws.Range("N94").Value = -2887
$ws.Range("H103").Value = 22000
$ws.Range("J103").Value = 22000
$ws.Range("L103").Value = 22000
$ws.Range("N103").Value = -24344
$ws.Range("H112").Value = 29533
$ws.Range("J112").Value = 29533
$ws.Range("L112").Value = 29533
$ws.Range("N112").Value = -32487
$ws.Range("H118").Value = 47712
$ws.Range("J118").Value = 47712
$ws.Range("L118").Value = 47712
$ws.Range("N118").Value = -51026
$ws.Range("H130").Value = 56513.332
$ws.Range("J130").Value = 56513.332
$ws.Range("L130").Value = 56513.332
$ws.Range("N130").Value = -66553.33199999999
$ws.Range("H132").Value = 71990
$ws.Range("J132").Value = 71990
$ws.Range("L132").Value = 71990
$ws.Range("N132").Value = -82110
$ws.Range("H135").Value = 57822.5
$ws.Range("J135").Value = 57822.5
$ws.Range("L135").Value = 57822.5
$ws.Range("N135").Value = -67962.5
$ws.Range("H137").Value = 64780
$ws.Range("J137").Value = 64780
$ws.Range("L137").Value = 64780
$ws.Range("N137").Value = -74980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = -10226
$ws.Range("H86").Value = 3001.75
$ws.Range("I86").Value = 3002.3333
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3002.3333
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1879.3333
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3001.75
$ws.Range("I89").Value = 3002.3333
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 15011.6665
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -9395.666499999999
$ws.Range("N89").Value = -26232
$ws.Range("H99").Value = 2133.8635
$ws.Range("I99").Value = 2433.8125
$ws.Range("K99").Value = 2433.8125
$ws.Range("M99").Value = -935.8125
$ws.Range("H126").Value = 2133.8635
$ws.Range("I126").Value = 2433.8125
$ws.Range("K126").Value = 7301.4375
$ws.Range("M126").Value = -4831.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2100
$ws.Range("J62").Value = 2100
$ws.Range("L62").Value = 6300
$ws.Range("N62").Value = -7672
$ws.Range("H64").Value = 6975
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 6975
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 20925
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -21465
$ws.Range("H65").Value = 2100
$ws.Range("J65").Value = 2100
$ws.Range("L65").Value = 18900
$ws.Range("N65").Value = -25764
$ws.Range("H67").Value = 6975
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 6975
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 20925
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -22797
$ws.Range("H68").Value = 1393.875
$ws.Range("J68").Value = 1307.2858
$ws.Range("L68").Value = 3921.8574
$ws.Range("N68").Value = -5543.857400000001
$ws.Range("H71").Value = 1393.875
$ws.Range("J71").Value = 1307.2858
$ws.Range("L71").Value = 11765.5722
$ws.Range("N71").Value = -19877.5722
$ws.Range("H122").Value = 677.6875
$ws.Range("I122").Value = 421.73914
$ws.Range("K122").Value = 3795.65226
$ws.Range("M122").Value = -1345.65226

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1821.9459
$ws.Range("I102").Value = 1786.1666
$ws.Range("K102").Value = 1786.1666
$ws.Range("M102").Value = -164.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3479.55
$ws.Range("I7").Value = 2661.6155
$ws.Range("J7").Value = 4998.5713
$ws.Range("K7").Value = 2661.6155
$ws.Range("L7").Value = 4998.5713
$ws.Range("M7").Value = -2549.6155
$ws.Range("N7").Value = -5222.5713
$ws.Range("H25").Value = 450000
$ws.Range("I25").Value = 450000
$ws.Range("K25").Value = 450000
$ws.Range("M25").Value = -449770
$ws.Range("H40").Value = 9547.5
$ws.Range("I40").Value = 11730
$ws.Range("K40").Value = 11730
$ws.Range("M40").Value = -11594
$ws.Range("H82").Value = 2747.7856
$ws.Range("I82").Value = 1965.3334
$ws.Range("J82").Value = 4156.2
$ws.Range("K82").Value = 1965.3334
$ws.Range("L82").Value = 4156.2
$ws.Range("M82").Value = -1604.3334
$ws.Range("N82").Value = -4878.2
$ws.Range("H85").Value = 2747.7856
$ws.Range("I85").Value = 1965.3334
$ws.Range("J85").Value = 4156.2
$ws.Range("K85").Value = 1965.3334
$ws.Range("L85").Value = 4156.2
$ws.Range("M85").Value = -717.3334
$ws.Range("N85").Value = -6652.2
$ws.Range("H100").Value = 4899.4
$ws.Range("I100").Value = 7836.5
$ws.Range("J100").Value = 1542.7142
$ws.Range("K100").Value = 7836.5
$ws.Range("L100").Value = 1542.7142
$ws.Range("M100").Value = -7295.5
$ws.Range("N100").Value = -2624.7142
$ws.Range("H126").Value = 3479.55
$ws.Range("I126").Value = 2661.6155
$ws.Range("J126").Value = 4998.5713
$ws.Range("K126").Value = 7984.8465
$ws.Range("L126").Value = 14995.7139
$ws.Range("M126").Value = -5514.8465
$ws.Range("N126").Value = -19935.7139
$ws.Range("H132").Value = 4762.3257
$ws.Range("I132").Value = 4507.564
$ws.Range("J132").Value = 7246.25
$ws.Range("K132").Value = 13522.692
$ws.Range("L132").Value = 21738.75
$ws.Range("M132").Value = -10992.692
$ws.Range("N132").Value = -26798.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 59806.668
$ws.Range("I26").Value = 49604
$ws.Range("J26").Value = 70009.336
$ws.Range("K26").Value = 49604
$ws.Range("L26").Value = 70009.336
$ws.Range("M26").Value = -49311
$ws.Range("N26").Value = -70595.336
